$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: zip() ---------------------------------------------------------
# (values are entered A, B, D, C to reproduce the original shared-string order)
$ws.Range("A17").Value = "zip(이터러블)"
$ws.Range("B17").Value = "여러 개의 반복 가능한(iterable) 객체를 병렬로 묶어주는 함수"
$ws.Range("D17").Value = "zip(*iterables)"
$ws.Range("C17").Value = "각 반복 가능한 객체의 동일한 위치에 있는 요소들을 묶어 튜플 형태로 반환`nlist1, list2 = zip(*pairs) 형태로 언패킹도 가능"

# --- Row 18: dict.fromkeys() ------------------------------------------------
# (values are entered A, B, D, C to reproduce the original shared-string order)
$ws.Range("A18").Value = "dict.fromkeys()"
$ws.Range("B18").Value = "지정한 키들로 딕셔너리를 생성하고, 모든 키의 값을 동일한 값으로 설정하는 클래스 메서드"
$ws.Range("D18").Value = "dict.fromkeys(keys, value=None)"
$ws.Range("C18").Value = "기본 값이 동일한 키-값 쌍을 가진 딕셔너리 생성 가능"

# --- Row 19: dict.get() ------------------------------------------------------
$ws.Range("A19").Value = "dict.get()"
$ws.Range("B19").Value = "딕셔너리에서 지정한 키의 값을 가져오는 데 사용되는 메서드"
$ws.Range("C19").Value = "지정한 키가 딕셔너리에 존재하지 않을 경우, 기본적으로 None을 반환하며, 사용자 지정 기본값을 설정 가능"
$ws.Range("D19").Value = "dict.get(key, default=None)"

# Match font + layout of the existing table rows (A:D use the "D2Coding"
# bodied style - xf 1 - everywhere, except column C on row 17 which wraps
# the two-line description, same as the other multi-line example cells).
$ws.Range("A17:D19").Font.Name = "D2Coding"
$ws.Range("C17").WrapText = $true

# Row heights match the rest of the table (40pt, explicit custom height).
$ws.Rows.Item(17).RowHeight = 40
$ws.Rows.Item(18).RowHeight = 40
$ws.Rows.Item(19).RowHeight = 40

# Move the selection down to the new last cell, like the author would have
# left it after typing in the new rows.
$ws.Range("D20").Select()
